# Applies the row 6 / row 7 data swap described in the diff.
# Row 6 takes on the values that used to be in row 7, and vice versa,
# for the columns that actually differ between the two rows
# (A, B, E, F, G, H, Q, R). All other populated columns already hold
# identical values on both rows, so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for row 6 (previously held by row 7)
$ws.Range("A6").Value = 130905540
$ws.Range("B6").Value = 80348
$ws.Range("E6").Value = 6458
$ws.Range("F6").Value = "Lunglav"
$ws.Range("G6").Value = "Lobaria pulmonaria"
$ws.Range("H6").Value = "(L.) Hoffm."
$ws.Range("Q6").Value = 397970
$ws.Range("R6").Value = 7048600

# New values for row 7 (previously held by row 6)
$ws.Range("A7").Value = 130904210
$ws.Range("B7").Value = 79243
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("Q7").Value = 397838
$ws.Range("R7").Value = 7048574
